$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2-dose series")

# --- Update Series-level metadata -----------------------------------------
# Series renamed from the universal "primary series" to a Risk-based,
# high-priority 2-dose series (WHO SAGE prioritization: healthcare workers,
# 60+, immunocompromised).

# --- Insert rows for the new "Indication" entries --------------------------
# 3 new rows are needed right after "Select Patient Series" (row 6) for the
# Healthcare worker / Older adult 60+ / Immunocompromised indications, which
# pushes everything below down by 3 rows.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# A couple of rows shrink in column extent after the shift, so clear the
# leftover trailing cells carried down from the old layout:
#  - the "Age" row (old row 8, now row 11) only keeps columns A:F
#  - the former "Preferable Interval" row (old row 17, now row 20) only keeps columns A:I
$ws.Range("G11:H11").ClearContents()
$ws.Range("J20:L20").ClearContents()

# --- Write the final cell values -------------------------------------------
$ws.Cells.Item(1, 1).Value = "Series Name"
$ws.Cells.Item(1, 2).Value = "WHO COVID-19 high-priority 2-dose series"
$ws.Cells.Item(2, 1).Value = "Target Disease"
$ws.Cells.Item(2, 2).Value = "COVID-19"
$ws.Cells.Item(3, 1).Value = "Vaccine Group"
$ws.Cells.Item(3, 2).Value = "COVID-19"
$ws.Cells.Item(4, 1).Value = "Series Type"
$ws.Cells.Item(4, 2).Value = "Risk"
$ws.Cells.Item(5, 1).Value = "Equivalent Series Groups"
$ws.Cells.Item(5, 2).Value = "n/a"
$ws.Cells.Item(6, 1).Value = "Select Patient Series"
$ws.Cells.Item(6, 2).Value = "No"
$ws.Cells.Item(6, 3).Value = "No"
$ws.Cells.Item(6, 4).Value = "High Priority"
$ws.Cells.Item(6, 5).Formula = "=""1"""
$ws.Cells.Item(6, 5).Copy()
$ws.Cells.Item(6, 5).PasteSpecial(-4163)
$ws.Cells.Item(6, 6).Value = "A"
$ws.Cells.Item(6, 7).Formula = "=""1"""
$ws.Cells.Item(6, 7).Copy()
$ws.Cells.Item(6, 7).PasteSpecial(-4163)
$ws.Cells.Item(6, 8).Value = "n/a"
$ws.Cells.Item(6, 9).Value = "n/a"
$ws.Cells.Item(7, 1).Value = "Indication"
$ws.Cells.Item(7, 2).Value = "Healthcare worker (1020)"
$ws.Cells.Item(7, 3).Value = "Patient is a healthcare worker with occupational exposure risk"
$ws.Cells.Item(7, 4).Value = "n/a"
$ws.Cells.Item(7, 5).Value = "n/a"
$ws.Cells.Item(8, 1).Value = "Indication"
$ws.Cells.Item(8, 2).Value = "Older adult 60+ years (1021)"
$ws.Cells.Item(8, 3).Value = "Patient is an older adult (60 years or older)"
$ws.Cells.Item(8, 4).Value = "60 years"
$ws.Cells.Item(8, 5).Value = "n/a"
$ws.Cells.Item(9, 1).Value = "Indication"
$ws.Cells.Item(9, 2).Value = "Immunocompromised individual (1022)"
$ws.Cells.Item(9, 3).Value = "Patient is immunocompromised"
$ws.Cells.Item(9, 4).Value = "n/a"
$ws.Cells.Item(9, 5).Value = "n/a"
$ws.Cells.Item(10, 1).Value = "Series Dose"
$ws.Cells.Item(10, 2).Value = "Dose 1"
$ws.Cells.Item(11, 1).Value = "Age"
$ws.Cells.Item(11, 2).Value = "12 years"
$ws.Cells.Item(11, 3).Value = "18 years"
$ws.Cells.Item(11, 4).Value = "18 years"
$ws.Cells.Item(11, 5).Value = "n/a"
$ws.Cells.Item(11, 6).Value = "n/a"
$ws.Cells.Item(12, 1).Value = "Preferable Vaccine"
$ws.Cells.Item(12, 2).Value = "COVID-19, mRNA (Pfizer) (208)"
$ws.Cells.Item(12, 3).Value = "12 years"
$ws.Cells.Item(12, 4).Value = "n/a"
$ws.Cells.Item(12, 5).Value = "n/a"
$ws.Cells.Item(12, 6).Value = "n/a"
$ws.Cells.Item(12, 7).Value = "Y"
$ws.Cells.Item(13, 1).Value = "Preferable Vaccine"
$ws.Cells.Item(13, 2).Value = "COVID-19, mRNA (Moderna) (207)"
$ws.Cells.Item(13, 3).Value = "12 years"
$ws.Cells.Item(13, 4).Value = "n/a"
$ws.Cells.Item(13, 5).Value = "n/a"
$ws.Cells.Item(13, 6).Value = "n/a"
$ws.Cells.Item(13, 7).Value = "N"
$ws.Cells.Item(14, 1).Value = "Allowable Vaccine"
$ws.Cells.Item(14, 2).Value = "COVID-19, mRNA (Pfizer) (208)"
$ws.Cells.Item(14, 3).Value = "12 years"
$ws.Cells.Item(14, 4).Value = "n/a"
$ws.Cells.Item(15, 1).Value = "Allowable Vaccine"
$ws.Cells.Item(15, 2).Value = "COVID-19, mRNA (Moderna) (207)"
$ws.Cells.Item(15, 3).Value = "12 years"
$ws.Cells.Item(15, 4).Value = "n/a"
$ws.Cells.Item(16, 1).Value = "Allowable Vaccine"
$ws.Cells.Item(16, 2).Value = "COVID-19, viral vector (J&J) (212)"
$ws.Cells.Item(16, 3).Value = "18 years"
$ws.Cells.Item(16, 4).Value = "n/a"
$ws.Cells.Item(17, 1).Value = "Allowable Vaccine"
$ws.Cells.Item(17, 2).Value = "COVID-19, protein subunit (Novavax) (211)"
$ws.Cells.Item(17, 3).Value = "12 years"
$ws.Cells.Item(17, 4).Value = "n/a"
$ws.Cells.Item(18, 1).Value = "Recurring Dose"
$ws.Cells.Item(18, 2).Value = "No"
$ws.Cells.Item(19, 1).Value = "Series Dose"
$ws.Cells.Item(19, 2).Value = "Dose 2"
$ws.Cells.Item(20, 1).Value = "Preferable Interval"
$ws.Cells.Item(20, 2).Value = "Y"
$ws.Cells.Item(20, 3).Value = "n/a"
$ws.Cells.Item(20, 4).Value = "n/a"
$ws.Cells.Item(20, 5).Value = "n/a"
$ws.Cells.Item(20, 6).Value = "3 weeks"
$ws.Cells.Item(20, 7).Value = "4 weeks"
$ws.Cells.Item(20, 8).Value = "4 weeks"
$ws.Cells.Item(20, 9).Value = "n/a"
$ws.Cells.Item(21, 1).Value = "Preferable Vaccine"
$ws.Cells.Item(21, 2).Value = "COVID-19, mRNA (Pfizer) (208)"
$ws.Cells.Item(21, 3).Value = "12 years"
$ws.Cells.Item(21, 4).Value = "n/a"
$ws.Cells.Item(21, 5).Value = "n/a"
$ws.Cells.Item(21, 6).Value = "n/a"
$ws.Cells.Item(21, 7).Value = "Y"
$ws.Cells.Item(22, 1).Value = "Preferable Vaccine"
$ws.Cells.Item(22, 2).Value = "COVID-19, mRNA (Moderna) (207)"
$ws.Cells.Item(22, 3).Value = "12 years"
$ws.Cells.Item(22, 4).Value = "n/a"
$ws.Cells.Item(22, 5).Value = "n/a"
$ws.Cells.Item(22, 6).Value = "n/a"
$ws.Cells.Item(22, 7).Value = "N"
$ws.Cells.Item(23, 1).Value = "Allowable Vaccine"
$ws.Cells.Item(23, 2).Value = "COVID-19, mRNA (Pfizer) (208)"
$ws.Cells.Item(23, 3).Value = "12 years"
$ws.Cells.Item(23, 4).Value = "n/a"
$ws.Cells.Item(24, 1).Value = "Allowable Vaccine"
$ws.Cells.Item(24, 2).Value = "COVID-19, mRNA (Moderna) (207)"
$ws.Cells.Item(24, 3).Value = "12 years"
$ws.Cells.Item(24, 4).Value = "n/a"
$ws.Cells.Item(25, 1).Value = "Allowable Vaccine"
$ws.Cells.Item(25, 2).Value = "COVID-19, protein subunit (Novavax) (211)"
$ws.Cells.Item(25, 3).Value = "12 years"
$ws.Cells.Item(25, 4).Value = "n/a"
$ws.Cells.Item(26, 1).Value = "Recurring Dose"
$ws.Cells.Item(26, 2).Value = "No"
